# Fruta / hortaliza, semanal
# Insert two new daily price rows for Piña (Femacal de La Calera) ahead of
# the existing row 1107, shifting the remaining rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 1107:1200 down to 1109:1202, leaving two blank rows.
$ws.Rows("1107:1108").Insert()

# New row 1107
$ws.Cells.Item(1107, 1).Value = 3
$ws.Cells.Item(1107, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1107, 3).Value = "Coquimbo"
$ws.Cells.Item(1107, 4).Value = 45223
$ws.Cells.Item(1107, 5).Value = 5
$ws.Cells.Item(1107, 6).Value = "Fruta"
$ws.Cells.Item(1107, 7).Value = 100108
$ws.Cells.Item(1107, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(1107, 9).Value = 100108005
$ws.Cells.Item(1107, 10).Value = "Piña"
$ws.Cells.Item(1107, 11).Value = "Caramelo"
$ws.Cells.Item(1107, 12).Value = "Especial"
$ws.Cells.Item(1107, 13).Value = 80
$ws.Cells.Item(1107, 14).Value = 21000
$ws.Cells.Item(1107, 15).Value = 21000
$ws.Cells.Item(1107, 16).Value = 21000
$ws.Cells.Item(1107, 17).Value = "$/caja 10 unidades"
$ws.Cells.Item(1107, 18).Value = "Ecuador"
$ws.Cells.Item(1107, 19).Value = 2100
$ws.Cells.Item(1107, 20).Value = 10

# New row 1108
$ws.Cells.Item(1108, 1).Value = 3
$ws.Cells.Item(1108, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1108, 3).Value = "Coquimbo"
$ws.Cells.Item(1108, 4).Value = 45223
$ws.Cells.Item(1108, 5).Value = 5
$ws.Cells.Item(1108, 6).Value = "Fruta"
$ws.Cells.Item(1108, 7).Value = 100108
$ws.Cells.Item(1108, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(1108, 9).Value = 100108005
$ws.Cells.Item(1108, 10).Value = "Piña"
$ws.Cells.Item(1108, 11).Value = "Caramelo"
$ws.Cells.Item(1108, 12).Value = "Primera"
$ws.Cells.Item(1108, 13).Value = 60
$ws.Cells.Item(1108, 14).Value = 21000
$ws.Cells.Item(1108, 15).Value = 21000
$ws.Cells.Item(1108, 16).Value = 21000
$ws.Cells.Item(1108, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(1108, 18).Value = "Ecuador"
$ws.Cells.Item(1108, 19).Value = 1750
$ws.Cells.Item(1108, 20).Value = 12
